# Update the "releve" bank-statement export with the new statement data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal TEXT, never letting Excel's automatic
# type inference turn date-looking / number-looking strings (e.g.
# "02/01/2025", "0823752100109", "31") into real dates or numbers.
function Set-TextValue {
    param(
        [string]$Addr,
        [string]$Val
    )
    $r = $ws.Range($Addr)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.Style = "Normal"
}

# --- Header block ---------------------------------------------------------
Set-TextValue "B1" "AFRILAND FIRST BANK"
Set-TextValue "B2" "0823752100109"
Set-TextValue "B3" "SAFIR CONSULTING CAMEROUN"
Set-TextValue "B4" "02/03/2025 - 02/01/2025"

# --- Table header (row 7) stays Date / Description / Montant / Sens ------
Set-TextValue "A7" "Date"
Set-TextValue "B7" "Description"
Set-TextValue "C7" "Montant"
Set-TextValue "D7" "Sens"

# --- Transaction rows (8-20) -----------------------------------------------
$rows = @(
    @{ Row = 8;  A = "02/01/2025"; B = "/24 PRELVALIOS FINANCE 12 1257225";            C = "31" }
    @{ Row = 9;  A = "02/01/2025"; B = "/24 VIREMENT SALAIRE. 3 225 MENSUEL SAFI";      C = "31" }
    @{ Row = 10; A = "03/09/2025"; B = "/24 FRAIS FIXE AU /24 5000 952.225";            C = "31" }
    @{ Row = 11; A = "03/01/2025"; B = "/24 COMMISSION DE CPTE AU 107 BunzR4";          C = "31" }
    @{ Row = 12; A = "03/01/2025"; B = "/24 COMM. DE DECGUVERT AU 40 » 31224";          C = "31" }
    @{ Row = 13; A = "31/12/2024"; B = "0370172025 INTERETS DEBITEURS AU ( 909 771";    C = "28163" }
    @{ Row = 14; A = "12/24";      B = "/24";                                           C = "3112" }
    @{ Row = 15; A = "03/01/2025"; B = "/24 TAXE/INTERETS OBT AU 54";                   C = "31" }
    @{ Row = 16; A = "12/24";      B = "/24";                                           C = "3112" }
    @{ Row = 17; A = "03/01/2028"; B = "«/24 TX/COM. DECOUVERT Ais 7 31224";            C = "31" }
    @{ Row = 18; A = "03/01/2025"; B = "/24 TAXE COMM. DE CPTE AU 19 Bina";             C = "31" }
    @{ Row = 19; A = "03/01/2025"; B = "/24 TAXE FRAIS FIXE AU 9 B12";                  C = "31" }
    @{ Row = 20; A = "06/01/2025"; B = "/25 ViREMENT CIME BONA 11925, BBR 7IT";         C = "03" }
)

foreach ($row in $rows) {
    Set-TextValue ("A" + $row.Row) $row.A
    Set-TextValue ("B" + $row.Row) $row.B
    Set-TextValue ("C" + $row.Row) $row.C
}

# Clear the now-unused "Sens" (D) column for every transaction row - the
# new statement layout only uses Date / Description / Montant.
$ws.Range("D8:D11").ClearContents()
